$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: item names (filled first so shared-string order matches source)
$ws.Range("A1").Value = "ItemName"
$ws.Range("A2").Value = "Some name."
$ws.Range("A3").Value = "Some name."
$ws.Range("A4").Value = "Some name."
$ws.Range("A5").Value = "Some name."
$ws.Range("A6").Value = "Some name."

# Column B: item prices
$ws.Range("B1").Value = "ItemPrice"
$ws.Range("B2").Value = 12
$ws.Range("B3").Value = 12
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 12

# Currency formatting on price column
$ws.Range("B2:B6").Style = "Currency"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 21.14
$ws.Columns.Item(2).ColumnWidth = 14.666666

# Selection as recorded in the source workbook
$ws.Range("F5").Select() | Out-Null
